$wb = $excel.ActiveWorkbook

# ALC!row44
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 45000
$ws.Range("J44").Value = 45000
$ws.Range("L44").Value = 45000
$ws.Range("N44").Value = -45924

# ALC!row51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9400.4
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 10500.5
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 10500.5
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -11468.5

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5863.636
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 6750
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 6750
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -13258

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 721927.9
$ws.Range("I116").Value = 1668531.6
$ws.Range("J116").Value = 11975
$ws.Range("K116").Value = 1668531.6
$ws.Range("L116").Value = 11975
$ws.Range("M116").Value = -1665089.6
$ws.Range("N116").Value = -18859

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 879.2308
$ws.Range("I2").Value = 785.8333
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 785.8333
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -672.8333
$ws.Range("N2").Value = -2226

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2862.3333
$ws.Range("I45").Value = 3158.8
$ws.Range("J45").Value = 1380
$ws.Range("K45").Value = 3158.8
$ws.Range("L45").Value = 1380
$ws.Range("M45").Value = -2781.8
$ws.Range("N45").Value = -2134

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5589.0454
$ws.Range("I74").Value = 7714.9165
$ws.Range("J74").Value = 3038
$ws.Range("K74").Value = 7714.9165
$ws.Range("L74").Value = 3038
$ws.Range("M74").Value = -6840.9165
$ws.Range("N74").Value = -4786

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5589.0454
$ws.Range("I77").Value = 7714.9165
$ws.Range("J77").Value = 3038
$ws.Range("K77").Value = 38574.5825
$ws.Range("L77").Value = 15190
$ws.Range("M77").Value = -34206.5825
$ws.Range("N77").Value = -23926

# ARM!row101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 879.2308
$ws.Range("I116").Value = 785.8333
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 785.8333
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1508.1667
$ws.Range("N116").Value = -6588

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2372.875
$ws.Range("I132").Value = 1025.4
$ws.Range("J132").Value = 3335.3572
$ws.Range("K132").Value = 3076.2
$ws.Range("L132").Value = 10006.0716
$ws.Range("M132").Value = -546.2000000000003
$ws.Range("N132").Value = -15066.0716

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 879.2308
$ws.Range("I3").Value = 785.8333
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 785.8333
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -671.8333
$ws.Range("N3").Value = -2228

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1895.8718
$ws.Range("I134").Value = 1237.4642
$ws.Range("K134").Value = 3712.3926
$ws.Range("M134").Value = -1177.3926

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2748.5278
$ws.Range("I31").Value = 990.4737
$ws.Range("J31").Value = 4713.4116
$ws.Range("K31").Value = 990.4737
$ws.Range("L31").Value = 4713.4116
$ws.Range("M31").Value = -695.4737
$ws.Range("N31").Value = -5303.4116

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2748.5278
$ws.Range("I34").Value = 990.4737
$ws.Range("J34").Value = 4713.4116
$ws.Range("K34").Value = 990.4737
$ws.Range("L34").Value = 4713.4116
$ws.Range("M34").Value = -788.4737
$ws.Range("N34").Value = -5117.4116

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2499.527
$ws.Range("J58").Value = 5335.5
$ws.Range("L58").Value = 5335.5
$ws.Range("N58").Value = -5741.5

# CRP!row74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 40504.8
$ws.Range("J74").Value = 40504.8
$ws.Range("L74").Value = 40504.8
$ws.Range("N74").Value = -42252.8

# CRP!row77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 40504.8
$ws.Range("J77").Value = 40504.8
$ws.Range("L77").Value = 121514.4
$ws.Range("N77").Value = -130250.4

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1794.55
$ws.Range("I105").Value = 1431.3125
$ws.Range("K105").Value = 1431.3125
$ws.Range("M105").Value = 315.6875

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2325.3044
$ws.Range("I132").Value = 1324.2
$ws.Range("K132").Value = 3972.6
$ws.Range("M132").Value = -1442.6

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2499.527
$ws.Range("J136").Value = 5335.5
$ws.Range("L136").Value = 16006.5
$ws.Range("N136").Value = -21106.5

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 111
$ws.Range("I12").Value = 43.666668
$ws.Range("K12").Value = 131.000004
$ws.Range("M12").Value = 41.99999600000001

# CUL!row50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 282.66666
$ws.Range("I50").Value = 298.8
$ws.Range("J50").Value = 262.5
$ws.Range("K50").Value = 896.4000000000001
$ws.Range("L50").Value = 787.5
$ws.Range("M50").Value = -415.4000000000001
$ws.Range("N50").Value = -1749.5

# CUL!row53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 282.66666
$ws.Range("I53").Value = 298.8
$ws.Range("J53").Value = 262.5
$ws.Range("K53").Value = 896.4000000000001
$ws.Range("L53").Value = 787.5
$ws.Range("M53").Value = -415.4000000000001
$ws.Range("N53").Value = -1749.5

# CUL!row57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1164.5714
$ws.Range("I57").Value = 275
$ws.Range("J57").Value = 1520.4
$ws.Range("K57").Value = 825
$ws.Range("L57").Value = 4561.200000000001
$ws.Range("M57").Value = -266
$ws.Range("N57").Value = -5679.200000000001

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3906849.8
$ws.Range("I113").Value = 616.7059
$ws.Range("K113").Value = 1850.1177
$ws.Range("M113").Value = 319.8822999999998

# CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1762.3182
$ws.Range("I139").Value = 1110
$ws.Range("J139").Value = 2903.875
$ws.Range("K139").Value = 3330
$ws.Range("L139").Value = 8711.625
$ws.Range("M139").Value = 1810
$ws.Range("N139").Value = -18991.625

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3060.8928
$ws.Range("I132").Value = 1788.6471
$ws.Range("K132").Value = 5365.9413
$ws.Range("M132").Value = -2835.9413

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5696.174
$ws.Range("I132").Value = 1210.3334
$ws.Range("K132").Value = 3631.0002
$ws.Range("M132").Value = -1101.0002

# LTW!row133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 30680.25
$ws.Range("J133").Value = 30680.25
$ws.Range("L133").Value = 30680.25
$ws.Range("N133").Value = -35740.25

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4061.1924
$ws.Range("I136").Value = 1129.75
$ws.Range("K136").Value = 3389.25
$ws.Range("M136").Value = -839.25

# WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 728.625
$ws.Range("I113").Value = 728.625
$ws.Range("K113").Value = 2185.875
$ws.Range("M113").Value = -15.875

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6805273
$ws.Range("I132").Value = 1723.5927
$ws.Range("J132").Value = 15155084
$ws.Range("K132").Value = 5170.7781
$ws.Range("L132").Value = 45465252
$ws.Range("M132").Value = -2640.7781
$ws.Range("N132").Value = -45470312

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3974.6428
$ws.Range("I136").Value = 2108.7058
$ws.Range("K136").Value = 6326.117400000001
$ws.Range("M136").Value = -3776.117400000001
